$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update rows 2-34 in place, then drop old rows 35-37 ---
$ws1.Range("A2").Value = "BRVM-PRINCIPAL     (**)"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 674.21
$ws1.Range("E2").Value = 227.3
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 673.88
$ws1.Range("E3").Value = 172.07
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 672.0700000000001
$ws1.Range("E4").Value = 228.72
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 590.64
$ws1.Range("E5").Value = 148.16
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "BRVM-PRESTIGE"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 573.4
$ws1.Range("E6").Value = 144.05
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 572.5
$ws1.Range("E7").Value = 145.9
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "BRVM - ENERGIE"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 454.83
$ws1.Range("E8").Value = 113.19
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 449.02
$ws1.Range("E9").Value = 113.3
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

$ws1.Range("A10").Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 3
$ws1.Range("D10").Value = 402.69
$ws1.Range("E10").Value = 135.14
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 374.47
$ws1.Range("E11").Value = 94.3
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

$ws1.Range("A12").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B12").Value = 3
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 21.82
$ws1.Range("E12").Value = 6.83
$ws1.Range("F12").Value = "🟢 Achat"
$ws1.Range("G12").Value = "✅ Renforcer"

$ws1.Range("A13").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Range("B13").Value = 3
$ws1.Range("C13").Value = 1
$ws1.Range("D13").Value = 14.97
$ws1.Range("E13").Value = 7.38
$ws1.Range("F13").Value = "🟢 Achat"
$ws1.Range("G13").Value = "✅ Renforcer"

$ws1.Range("A14").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B14").Value = 2
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 8.41
$ws1.Range("E14").Value = 4.68
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

$ws1.Range("A15").Value = "SICOR CI (SICC)"
$ws1.Range("B15").Value = 2
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 7.19
$ws1.Range("E15").Value = 3.11
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 6.19
$ws1.Range("E16").Value = 6.19
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

$ws1.Range("A17").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 5.25
$ws1.Range("E17").Value = 5.25
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

$ws1.Range("A18").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B18").Value = 2
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = 4.75
$ws1.Range("E18").Value = 4.55
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "👀 À surveiller"

$ws1.Range("A19").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B19").Value = 1
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 3.15
$ws1.Range("E19").Value = -2.81
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "👀 À surveiller"

$ws1.Range("A20").Value = "NESTLE CI (NTLC)"
$ws1.Range("B20").Value = 2
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 2.27
$ws1.Range("E20").Value = 3.64
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "👀 À surveiller"

$ws1.Range("A21").Value = "SETAO CI (STAC)"
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 1.07
$ws1.Range("E21").Value = -4.35
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "👀 À surveiller"

$ws1.Range("A22").Value = "SICABLE CI (CABC)"
$ws1.Range("B22").Value = 2
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = -0.68
$ws1.Range("E22").Value = -3.35
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "👀 À surveiller"

$ws1.Range("A23").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = -0.7
$ws1.Range("E23").Value = -0.7
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"

$ws1.Range("A24").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = -1.43
$ws1.Range("E24").Value = -1.43
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"

$ws1.Range("A25").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B25").Value = 1
$ws1.Range("C25").Value = 2
$ws1.Range("D25").Value = -1.63
$ws1.Range("E25").Value = 6.08
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "👀 À surveiller"

$ws1.Range("A26").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B26").Value = 0
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = -1.88
$ws1.Range("E26").Value = -1.88
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"

$ws1.Range("A27").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B27").Value = 0
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = -1.92
$ws1.Range("E27").Value = -1.92
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"

$ws1.Range("A28").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = -2.44
$ws1.Range("E28").Value = -2.44
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "BERNABE CI (BNBC)"
$ws1.Range("B29").Value = 0
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = -2.53
$ws1.Range("E29").Value = -2.53
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

$ws1.Range("A30").Value = "SMB CI (SMBC)"
$ws1.Range("B30").Value = 0
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = -2.57
$ws1.Range("E30").Value = -2.57
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"

$ws1.Range("A31").Value = "SAFCA CI (SAFC)"
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -2.69
$ws1.Range("E31").Value = -2.69
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "➖ Neutre"

$ws1.Range("A32").Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = -3.49
$ws1.Range("E32").Value = -3.49
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"

$ws1.Range("A33").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Range("B33").Value = 0
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = -3.89
$ws1.Range("E33").Value = -3.89
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "➖ Neutre"

$ws1.Range("A34").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 2
$ws1.Range("D34").Value = -8.43
$ws1.Range("E34").Value = -4.26
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "➖ Neutre"

# Remove now-obsolete trailing rows (old rows 35, 36, 37)
$ws1.Rows(37).Delete()
$ws1.Rows(36).Delete()
$ws1.Rows(35).Delete()

# --- Sheet "Top_YTD": update rows 2-11 in place ---
$ws2.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B2").Value = 5093.44

$ws2.Range("A3").Value = "BRVM - SERVICES FINANCIERS"
$ws2.Range("B3").Value = 3662

$ws2.Range("A4").Value = "BRVM-PRESTIGE"
$ws2.Range("B4").Value = 3406.89

$ws2.Range("A5").Value = "BRVM - INDUSTRIELS"
$ws2.Range("B5").Value = 3393.05

$ws2.Range("A6").Value = "BRVM-PRINCIPAL     (**)"
$ws2.Range("B6").Value = 3324.26

$ws2.Range("A7").Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws2.Range("B7").Value = 3301.4

$ws2.Range("A8").Value = "BRVM - ENERGIE"
$ws2.Range("B8").Value = 1985.82

$ws2.Range("A9").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B9").Value = 1929.57

$ws2.Range("A10").Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Range("B10").Value = 1305.3

$ws2.Range("A11").Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws2.Range("B11").Value = 1185.05

